$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2981.9285
$ws.Range("I70").Value = 3203.6667
$ws.Range("J70").Value = 1651.5
$ws.Range("K70").Value = 9611.000100000001
$ws.Range("L70").Value = 4954.5
$ws.Range("M70").Value = -9341.000100000001
$ws.Range("N70").Value = -5494.5

$ws.Range("H73").Value = 2981.9285
$ws.Range("I73").Value = 3203.6667
$ws.Range("J73").Value = 1651.5
$ws.Range("K73").Value = 9611.000100000001
$ws.Range("L73").Value = 4954.5
$ws.Range("M73").Value = -8675.000100000001
$ws.Range("N73").Value = -6826.5

$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988

$ws.Range("H138").Value = 1885.4
$ws.Range("I138").Value = 801.4545
$ws.Range("J138").Value = 4156.524
$ws.Range("K138").Value = 2404.3635
$ws.Range("L138").Value = 12469.572
$ws.Range("M138").Value = 2735.6365
$ws.Range("N138").Value = -22749.572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H60").Value = 25557
$ws.Range("J60").Value = 25557
$ws.Range("L60").Value = 25557
$ws.Range("N60").Value = -27023

$ws.Range("H74").Value = 2554.1667
$ws.Range("I74").Value = 845.2549
$ws.Range("J74").Value = 12238
$ws.Range("K74").Value = 845.2549
$ws.Range("L74").Value = 12238
$ws.Range("M74").Value = 28.74509999999998
$ws.Range("N74").Value = -13986

$ws.Range("H77").Value = 2554.1667
$ws.Range("I77").Value = 845.2549
$ws.Range("J77").Value = 12238
$ws.Range("K77").Value = 4226.2745
$ws.Range("L77").Value = 61190
$ws.Range("M77").Value = 141.7254999999996
$ws.Range("N77").Value = -69926

$ws.Range("H104").Value = 37612.5
$ws.Range("J104").Value = 37612.5
$ws.Range("L104").Value = 37612.5
$ws.Range("N104").Value = -44600.5

$ws.Range("H132").Value = 17884.77
$ws.Range("I132").Value = 17553.715
$ws.Range("J132").Value = 18271
$ws.Range("K132").Value = 52661.145
$ws.Range("L132").Value = 54813
$ws.Range("M132").Value = -50131.145
$ws.Range("N132").Value = -59873

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 35403.516
$ws.Range("I20").Value = 894.1111
$ws.Range("J20").Value = 91873.45
$ws.Range("K20").Value = 894.1111
$ws.Range("L20").Value = 91873.45
$ws.Range("M20").Value = -647.1111
$ws.Range("N20").Value = -92367.45

$ws.Range("H106").Value = 21289.111
$ws.Range("J106").Value = 21289.111
$ws.Range("L106").Value = 21289.111
$ws.Range("N106").Value = -23813.111

$ws.Range("H107").Value = 2184.353
$ws.Range("I107").Value = 1050.375
$ws.Range("J107").Value = 3192.3333
$ws.Range("K107").Value = 1050.375
$ws.Range("L107").Value = 3192.3333
$ws.Range("M107").Value = 869.625
$ws.Range("N107").Value = -7032.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1305.6809
$ws.Range("I58").Value = 1084.3462
$ws.Range("J58").Value = 1579.7142
$ws.Range("K58").Value = 1084.3462
$ws.Range("L58").Value = 1579.7142
$ws.Range("M58").Value = -881.3462
$ws.Range("N58").Value = -1985.7142

$ws.Range("H70").Value = 11882.353
$ws.Range("J70").Value = 11882.353
$ws.Range("L70").Value = 11882.353
$ws.Range("N70").Value = -12512.353

$ws.Range("H73").Value = 11882.353
$ws.Range("J73").Value = 11882.353
$ws.Range("L73").Value = 11882.353
$ws.Range("N73").Value = -14066.353

$ws.Range("H106").Value = 18601.625
$ws.Range("J106").Value = 18601.625
$ws.Range("L106").Value = 18601.625
$ws.Range("N106").Value = -21125.625

$ws.Range("H132").Value = 22226276
$ws.Range("I132").Value = 28576102
$ws.Range("J132").Value = 1878.8
$ws.Range("K132").Value = 85728306
$ws.Range("L132").Value = 5636.4
$ws.Range("M132").Value = -85725776
$ws.Range("N132").Value = -10696.4

$ws.Range("H136").Value = 1305.6809
$ws.Range("I136").Value = 1084.3462
$ws.Range("J136").Value = 1579.7142
$ws.Range("K136").Value = 3253.0386
$ws.Range("L136").Value = 4739.142599999999
$ws.Range("M136").Value = -703.0385999999999
$ws.Range("N136").Value = -9839.1426

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 1741.5
$ws.Range("I99").Value = 1741.5
$ws.Range("K99").Value = 5224.5
$ws.Range("M99").Value = -2978.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 36000
$ws.Range("J104").Value = 36000
$ws.Range("L104").Value = 36000
$ws.Range("N104").Value = -42988

$ws.Range("H105").Value = 48000
$ws.Range("J105").Value = 48000
$ws.Range("L105").Value = 48000
$ws.Range("N105").Value = -54988

$ws.Range("H132").Value = 4865.811
$ws.Range("I132").Value = 6707.6523
$ws.Range("J132").Value = 1839.9286
$ws.Range("K132").Value = 20122.9569
$ws.Range("L132").Value = 5519.7858
$ws.Range("M132").Value = -17592.9569
$ws.Range("N132").Value = -10579.7858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4020.34
$ws.Range("I132").Value = 3990.9048
$ws.Range("J132").Value = 4174.875
$ws.Range("K132").Value = 11972.7144
$ws.Range("L132").Value = 12524.625
$ws.Range("M132").Value = -9442.714399999999
$ws.Range("N132").Value = -17584.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2983.3333
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 2966.6667
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 2966.6667
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -4214.6667

$ws.Range("H65").Value = 2983.3333
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 2966.6667
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 14833.3335
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -21073.3335

$ws.Range("H70").Value = 20000
$ws.Range("J70").Value = 20000
$ws.Range("L70").Value = 20000
$ws.Range("N70").Value = -20630

$ws.Range("H73").Value = 20000
$ws.Range("J73").Value = 20000
$ws.Range("L73").Value = 20000
$ws.Range("N73").Value = -22184

$ws.Range("H75").Value = 48000
$ws.Range("J75").Value = 48000
$ws.Range("L75").Value = 48000
$ws.Range("N75").Value = -49872

$ws.Range("H78").Value = 48000
$ws.Range("J78").Value = 48000
$ws.Range("L78").Value = 144000
$ws.Range("N78").Value = -153360

$ws.Range("H105").Value = 26123
$ws.Range("J105").Value = 26123
$ws.Range("L105").Value = 26123
$ws.Range("N105").Value = -33111
